# CWR_Checklist_Template.test.xlsx edit script
# - Adds two near-empty rows into "CK_Identification" (rows 10 and 12 in the
#   final layout), pushing the existing data rows down.
# - Rebuilds the worksheet's hyperlinks so they keep tracking the cells that
#   now hold the "Test dataset" / "http://bioversityinternational.org" rows.
# - Switches the active tab from "CK_Threats" to "CK_Crossability".
# - Restores the workbook window geometry recorded by the last save.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. CK_Identification ("Checklist identification" aka SessionUpload ws):
#    insert two rows into the small sample table (rows 10 & 12), trimming
#    the newly-inserted rows down to only the handful of cells that should
#    keep formatting but stay blank.
# ---------------------------------------------------------------------
$wsId = $wb.Worksheets.Item("CK_Identification")

# Insert a blank row at 10 -- old row10 (Test dataset / 201504) shifts to 11,
# old row11 (TEST-2 placeholder) shifts to 12, old row12 shifts to 13.
$wsId.Rows.Item(10).Insert()

# Insert a second blank row at 12 -- old row11 shifts again, from 12 to 13,
# old row12 shifts from 13 to 14.
$wsId.Rows.Item(12).Insert()

# New row 10 only keeps C10 / I10 / J10 / AE10 (all blank, formatted).
$wsId.Range("A10").Clear()
$wsId.Range("B10").Clear()
$wsId.Range("D10:H10").Clear()
$wsId.Range("K10:Z10").Clear()
$wsId.Range("AA10:AD10").Clear()
$wsId.Range("C10").ClearContents()
$wsId.Range("J10").ClearContents()

# New row 12 only keeps C12 / J12 (both blank, formatted).
$wsId.Range("B12").Clear()
$wsId.Range("D12:I12").Clear()
$wsId.Range("K12:Z12").Clear()
$wsId.Range("C12").ClearContents()
$wsId.Range("J12").ClearContents()

# Hyperlinks don't auto-track the row insert, so rebuild them against the
# cells that now hold the hyperlinked text.
$wsId.Hyperlinks.Delete()
$wsId.Hyperlinks.Add($wsId.Range("C9"), "http://pgrdiversity.bioversityinternational.org")
$wsId.Hyperlinks.Add($wsId.Range("I9"), "http://bioversityinternational.org")
$wsId.Hyperlinks.Add($wsId.Range("C11"), "http://pgrdiversity.bioversityinternational.org")
$wsId.Hyperlinks.Add($wsId.Range("C13"), "http://pgrdiversity.bioversityinternational.org")
$wsId.Hyperlinks.Add($wsId.Range("C14"), "http://pgrdiversity.bioversityinternational.org")
$wsId.Hyperlinks.Add($wsId.Range("I14"), "http://bioversityinternational.org")

# The frozen pane's next-free-row selection moves from A13 to A15.
$wsId.Range("A15").Select()

# ---------------------------------------------------------------------
# 2. Swap the active worksheet tab from CK_Threats to CK_Crossability.
# ---------------------------------------------------------------------
$wsCrossability = $wb.Worksheets.Item("CK_Crossability")
$wsCrossability.Activate()

# ---------------------------------------------------------------------
# 3. Restore the workbook window geometry saved with the workbook.
# ---------------------------------------------------------------------
$excel.Windows.Item(1).Left = 0
$excel.Windows.Item(1).Top = 0
$excel.Windows.Item(1).Width = 25600
$excel.Windows.Item(1).Height = 16060
